$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mimic the formatting of row 42 (the closest "plain" data row) for the new row 44
$ws.Range("A42:D42").Copy()
$ws.Range("A44:D44").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new row 44: 73. Set Matrix Zeroes / Java / 45026 (date)
$ws.Range("A44").Value = 73
$ws.Range("B44").Value = "Set Matrix Zeroes"
$ws.Range("C44").Value = "Java"
$ws.Range("D44").Value = 45026

# Update the selection shown in the saved workbook
$ws.Range("H39").Select()
